# Re-applies the refreshed "cryptos" price/volume snapshot (GitHub Actions bot).
# Only cell VALUES change (text content of columns B-E); no styling changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.558.17'
$ws.Range('E2').Value = '''  -0.94%  '
$ws.Range('D3').Value = '2.928.56'
$ws.Range('E3').Value = '''  -0.24%  '
$ws.Range('E4').Value = '''  -0.09%  '
$ws.Range('D5').Value = '''350.83'
$ws.Range('E5').Value = '''  -0.48%  '
$ws.Range('D6').Value = '''106.61'
$ws.Range('E6').Value = '''  -4.84%  '
$ws.Range('D7').Value = '''0.555'
$ws.Range('E7').Value = '''  -0.76%  '
$ws.Range('E8').Value = '''  -0.02%  '
$ws.Range('D9').Value = '''0.606'
$ws.Range('E9').Value = '''  -2.80%  '
$ws.Range('D10').Value = '''37.85'
$ws.Range('E10').Value = '''  -3.92%  '
$ws.Range('E11').Value = '''  +1.35%  '
$ws.Range('E12').Value = '''  -3.36%  '
$ws.Range('D13').Value = '''18.89'
$ws.Range('E13').Value = '''  -6.07%  '
$ws.Range('D14').Value = '3.380.70'
$ws.Range('E14').Value = '''  -0.69%  '
$ws.Range('D15').Value = '''7.53'
$ws.Range('E15').Value = '''  -2.87%  '
$ws.Range('D16').Value = '2.916.42'
$ws.Range('E16').Value = '''  -0.78%  '
$ws.Range('D17').Value = '''0.963'
$ws.Range('E17').Value = '''  -1.83%  '
$ws.Range('D18').Value = '51.512.05'
$ws.Range('E18').Value = '''  -1.23%  '
$ws.Range('D19').Value = '''3.38'
$ws.Range('E19').Value = '''  +2.67%  '
$ws.Range('E20').Value = '''  -2.86%  '
$ws.Range('D21').Value = '''13.43'
$ws.Range('E21').Value = '''  -5.92%  '
$ws.Range('D22').Value = '0.0₃0966'
$ws.Range('E22').Value = '''  -1.27%  '
$ws.Range('D23').Value = '''68.93'
$ws.Range('E23').Value = '''  -3.20%  '
$ws.Range('D24').Value = '''260.61'
$ws.Range('E24').Value = '''  -2.89%  '
$ws.Range('D25').Value = '''2.70'
$ws.Range('E25').Value = '''  -2.96%  '
$ws.Range('E26').Value = '''  -3.80%  '
$ws.Range('D27').Value = '''26.45'
$ws.Range('E27').Value = '''  -2.09%  '
$ws.Range('E28').Value = '''  +0.09%  '
$ws.Range('D29').Value = '''7.38'
$ws.Range('E29').Value = '''  +2.29%  '
$ws.Range('E30').Value = '''  -0.55%  '
$ws.Range('D31').Value = '''10.22'
$ws.Range('E31').Value = '''  -3.56%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = '''2.19'
$ws.Range('E32').Value = '''  -2.60%  '
$ws.Range('B33').Value = 'RenderToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D33').Value = '''6.05'
$ws.Range('E33').Value = '''  -1.15%  '
$ws.Range('D34').Value = '''35.70'
$ws.Range('E34').Value = '''  -3.45%  '
$ws.Range('D35').Value = '''50.41'
$ws.Range('E35').Value = '''  -4.87%  '
$ws.Range('D36').Value = '''0.0430'
$ws.Range('E36').Value = '''  -5.09%  '
$ws.Range('E37').Value = '''  -0.20%  '
$ws.Range('E38').Value = '''  -7.16%  '
$ws.Range('D39').Value = '''17.62'
$ws.Range('E39').Value = '''  -5.49%  '
$ws.Range('D40').Value = '''1.94'
$ws.Range('E40').Value = '''  -5.64%  '
$ws.Range('E41').Value = '''  -1.73%  '
$ws.Range('E42').Value = '''  -1.79%  '
$ws.Range('D43').Value = '''22.24'
$ws.Range('E43').Value = '''  -4.19%  '
$ws.Range('D44').Value = '''119.67'
$ws.Range('E44').Value = '''  +7.40%  '
$ws.Range('E45').Value = '''  -3.39%  '
$ws.Range('D46').Value = '2.094.49'
$ws.Range('E46').Value = '''  -4.70%  '
$ws.Range('E47').Value = '''  -5.89%  '
$ws.Range('E48').Value = '''  -9.35%  '
$ws.Range('E49').Value = '''  -3.46%  '
$ws.Range('D50').Value = '''0.0334'
$ws.Range('E50').Value = '''  -5.36%  '
$ws.Range('D51').Value = '''0.905'
$ws.Range('E51').Value = '''  -4.80%  '
